$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K data shifts to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (the old column D) into new D:E columns
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$cols = @("D","E","F","G","H","I","J","K","L","M")

$rowData = @{}
$rowData[7] = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
$rowData[8] = @(8601000, 8433000, 8125000, 8083000, 8239000, 7927000, 8157000, 8539000, 8549000, 8247000)
$rowData[9] = @(3564000, 3362000, 3071000, 3641000, 3406000, 3102000, 3254000, 3716000, 3910000, 3794000)
$rowData[10] = @(5037000, 5071000, 5054000, 4442000, 4833000, 4825000, 4903000, 4823000, 4639000, 4453000)
$rowData[11] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[12] = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
$rowData[13] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[14] = @(147000, 25000, 135000, 164000, -181000, 34000, -60000, 100000, 19000, -5000)
$rowData[15] = @(2370000, 2361000, 2330000, 2259000, 2173000, 2094000, 2054000, 2110000, 2092000, 1981000)
$rowData[16] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[17] = @(8122000, 7631000, 7310000, 7847000, 7512000, 7326000, 6994000, 8069000, 8238000, 7625000)
$rowData[18] = @(479000, 802000, 815000, 236000, 727000, 601000, 1163000, 470000, 311000, 622000)
$rowData[19] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[20] = @(32000, 55000, 42000, -9000, -42000, 44000, -52000, 27000, -60000, -15000)
$rowData[21] = @(2881000, 3218000, 3187000, 2486000, 2858000, 2739000, 3165000, 2607000, 2343000, 2588000)
$rowData[22] = @(664000, 633000, 637000, 576000, 581000, 595000, 613000, 631000, 619000, 630000)
$rowData[23] = @(-153000, 224000, 220000, -349000, 104000, 50000, 498000, -134000, -368000, -23000)
$rowData[24] = @(-8000, 17000, 47000, -412000, 38000, 98000, 292000, 149000, 111000, 119000)
$rowData[25] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[26] = @(-145000, 207000, 173000, 63000, 66000, -48000, 206000, -283000, -479000, -142000)
$rowData[27] = @(-141000, 196000, 176000, 69000, 72000, -48000, 206000, -283000, -479000, -142000)
$rowData[28] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[29] = @(0, 0, 0, 0, 7090000, "NA", "NA", "NA", "NA", "NA")
$rowData[30] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[31] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[32] = @(-32000, -55000, -42000, 9000, 42000, -44000, 52000, -27000, 60000, 15000)
$rowData[33] = @(-141000, 196000, 176000, 69000, 7162000, -48000, 206000, -283000, -479000, -142000)
$rowData[34] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[35] = @(-141000, 196000, 176000, 69000, 7162000, -48000, 206000, -283000, -479000, -142000)
$rowData[38] = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
$rowData[39] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[40] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[41] = @(6191000, 5726000, 4378000, 6610000, 4440000, 4802000, 2478000, 2870000, 3707000, 4006000)
$rowData[42] = @(632000, 3186000, 4008000, 2354000, 173000, 1610000, 4349000, 5444000, 2349000, 1650000)
$rowData[43] = @(3455000, 3555000, 3492000, 3711000, 3917000, 4118000, 4089000, 4138000, 1236000, 1004000)
$rowData[44] = @(919000, 859000, 622000, 1003000, 1009000, 751000, 979000, 1064000, 1296000, 981000)
$rowData[45] = @(1199000, 1121000, 895000, 575000, 626000, 654000, 601000, 601000, 1984000, 2215000)
$rowData[46] = @(12396000, 14447000, 13395000, 14253000, 10165000, 11935000, 12496000, 14117000, 10572000, 9856000)
$rowData[47] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[48] = @(21422000, 20816000, 20538000, 19925000, 19712000, 18901000, 18866000, 19209000, 19333000, 19176000)
$rowData[49] = @(49961000, 50046000, 50199000, 50360000, 50475000, 50498000, 50727000, 50484000, 50717000, 50977000)
$rowData[50] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[51] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[52] = @(2625000, 2542000, 2317000, 921000, 1067000, 1132000, 1235000, 1313000, 673000, 819000)
$rowData[53] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[54] = @(86404000, 87851000, 86449000, 85459000, 81419000, 82466000, 83324000, 85123000, 81295000, 80828000)
$rowData[55] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[56] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[57] = @(3637000, 4210000, 3143000, 3409000, 3176000, 2947000, 2616000, 3281000, 2894000, 2649000)
$rowData[58] = @(3596000, 5346000, 4846000, 3429000, 4036000, 4142000, 5125000, 5036000, 6554000, 7014000)
$rowData[59] = @(3467000, 3370000, 3658000, 3962000, 3859000, 3808000, 3830000, 4141000, 4189000, 4285000)
$rowData[60] = @(10700000, 12926000, 11647000, 10800000, 11071000, 10897000, 11571000, 12458000, 13637000, 13948000)
$rowData[61] = @(36288000, 35329000, 35771000, 37463000, 32825000, 34236000, 34459000, 35878000, 30759000, 29541000)
$rowData[62] = @(11087000, 11132000, 11086000, 10777000, 11218000, 18313000, 18279000, 17979000, 17903000, 17916000)
$rowData[63] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[64] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[65] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[66] = @(58134000, 59450000, 58556000, 59103000, 55184000, 63446000, 64309000, 66315000, 62299000, 61405000)
$rowData[67] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[68] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[69] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[70] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[71] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[72] = @(291000, 432000, 236000, -1255000, -1264000, -8426000, -8378000, -8584000, -8301000, -7822000)
$rowData[73] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[74] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[75] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[76] = @(28270000, 28401000, 27893000, 26356000, 26235000, 19020000, 19015000, 18808000, 18996000, 19423000)
$rowData[77] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[80] = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
$rowData[81] = @(-141000, 196000, 176000, 69000, 7162000, -48000, 206000, -283000, -479000, -142000)
$rowData[82] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[83] = @(2370000, 2361000, 2330000, 2259000, 2173000, 2094000, 2054000, 2110000, 2092000, 1981000)
$rowData[84] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[85] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[86] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[87] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[88] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[89] = @(2225000, 2927000, 2430000, 2653000, 2683000, 2802000, 1924000, 1268000, 650000, 1708000)
$rowData[90] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[91] = @(-1416000, -1266000, -1132000, -1532000, -682000, -608000, -497000, -7970000, -1259000, -845000)
$rowData[92] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[93] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[94] = @(-1143000, -1918000, -4369000, -7411000, 152000, 1562000, -438000, 3499000, -1836000, -1163000)
$rowData[95] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[96] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[97] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[98] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[99] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[100] = @(-570000, 339000, -295000, 3901000, -1680000, -1197000, -1234000, 1926000, 887000, -319000)
$rowData[101] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[102] = @(512000, 1348000, -2234000, 2147000, -362000, 2324000, -392000, -765000, -299000, 226000)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
